# Add a new "07-Apr" sheet, duplicated from the current last sheet ("03-Apr"),
# as the newest daily status-update tab, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# The most-recently-added day sheet is the last tab in the workbook ("03-Apr").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate it (same data/styles/columns) and place the copy right after it,
# mirroring how the original author rolled a new day forward each time.
$lastSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "07-Apr"

# The previously-active "03-Apr" tab is no longer selected; its whole grid
# becomes the stored selection instead of a single active cell.
$lastSheet.Cells.Select()

# The new "07-Apr" sheet becomes the active tab, with its whole grid selected.
$newSheet.Activate()
$newSheet.Cells.Select()
